$d = $word.ActiveDocument

# 1) "for arraignment on December 15, 2021." -> "... December 17, 2021." (single run)
$d.Content.Find.Execute("for arraignment on December 15, 2021.", $true, $false, $false, $false, $false, $true, 1, $false, "for arraignment on December 17, 2021.", 2)

# 2) Standalone "December 15, 2021" -> "December 17, 2021" (single run)
$d.Content.Find.Execute("December 15, 2021", $true, $false, $false, $false, $false, $true, 1, $false, "December 17, 2021", 2)

# 3) Judge's name: "Marianne" -> "Kyle", "Hemmeter" -> "Rohrer".
#    These two words sit in their own runs, sandwiched between plain single-space
#    runs that share identical run formatting ("Judge" <sp> "Marianne" <sp> "Hemmeter").
#    A plain text replace on this engine coalesces runs that end up with identical
#    formatting, which would merge "Judge"/" "/"Kyle"/" "/"Rohrer" into one run.
#    To keep the original run layout (one run per word/space, matching the source
#    document), replace the text first and then force each affected sub-range back
#    into its own run by toggling a character property (Bold on/off nets to the same
#    visible formatting but makes the engine re-split the run boundaries).

$rngMarianne = $d.Content
$rngMarianne.Find.Execute("Marianne", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngMarianne.Text = "Kyle"

$rngHemmeter = $d.Content
$rngHemmeter.Find.Execute("Hemmeter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngHemmeter.Text = "Rohrer"

# Re-locate each word/run after the text edits and force a clean run boundary.
$rngJudge = $d.Content
$rngJudge.Find.Execute("Judge", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngJudge.Font.Bold = $true
$rngJudge.Font.Bold = $false

$rngKyle = $d.Content
$rngKyle.Find.Execute("Kyle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngKyle.Font.Bold = $true
$rngKyle.Font.Bold = $false

$rngRohrer = $d.Content
$rngRohrer.Find.Execute("Rohrer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngRohrer.Font.Bold = $true
$rngRohrer.Font.Bold = $false

# Split the single-space runs that sit between Judge/Kyle and Kyle/Rohrer so they
# remain their own runs instead of bleeding into a neighboring word run.
$rngSpace1 = $d.Range($rngJudge.End, $rngKyle.Start)
$rngSpace1.Font.Bold = $true
$rngSpace1.Font.Bold = $false

$rngSpace2 = $d.Range($rngKyle.End, $rngRohrer.Start)
$rngSpace2.Font.Bold = $true
$rngSpace2.Font.Bold = $false
